$d = $word.ActiveDocument

# --- Locate the target paragraph ("One of the major failings ...") ---
$targetIndex = -1
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "One of the major failings*") {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)
$pStart = $p.Range.Start

# --- Re-split the first run: move the "_GoBack" bookmark so it sits right
#     after "...heat a" instead of at the end of the paragraph. Word
#     automatically breaks the run at the bookmark position on save. ---
$splitMarker = "One of the major failings of the equipment used to heat a"
$splitPos = $pStart + $splitMarker.Length
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Append the new trailing sentence to the end of the same paragraph
#     (just before its paragraph mark) ---
$p = $d.Paragraphs.Item($targetIndex)
$paraTextEnd = $p.Range.End - 1
$appendRange = $d.Range($paraTextEnd, $paraTextEnd)
$appendRange.InsertAfter(" This can lead to home-owners being disillusioned about their oil usage and can also lead to system problems such as air-locking in the houses plumbing eventually when the oil has run out. ")
$appendRange.Font.Name = "Times New Roman"
$appendRange.Font.NameBi = "Times New Roman"

# --- Insert a brand-new paragraph right after it ---
$p = $d.Paragraphs.Item($targetIndex)
$p.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$insPos = $newPara.Range.Start

# Run 1
$quoteOpen = [char]0x2018
$run1Text = "The idea for this project has been inspired by both a love for the $($quoteOpen)Internet "
$r1 = $d.Range($insPos, $insPos)
$r1.InsertAfter($run1Text)
$r1.Font.Name = "Times New Roman"
$r1.Font.NameBi = "Times New Roman"
$insPos = $insPos + $run1Text.Length

# Run 2 ("Of")
$r2 = $d.Range($insPos, $insPos)
$r2.InsertAfter("Of")
$r2.Font.Name = "Times New Roman"
$r2.Font.NameBi = "Times New Roman"
$insPos = $insPos + "Of".Length

# Run 3
$quoteClose = [char]0x2019
$r3 = $d.Range($insPos, $insPos)
$r3.InsertAfter(" Things$($quoteClose) and trying to support a greener environment, and through this project I am hoping I can accomplish both for people who use oil to heat their homes.")
$r3.Font.Name = "Times New Roman"
$r3.Font.NameBi = "Times New Roman"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
